$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by
# Excel's type inference (e.g. "0.999", "18.29") are first formatted as Text,
# written, then have that temporary formatting cleared again so the stored
# cell keeps its original (default/general) style while the value stays text.

$ws.Range("D2").Value = "68.525.89"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.778.31"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").Value = "3.776.84"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "4.412.79"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "3.781.29"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "68.449.16"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  +4.24%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  -4.31%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "3.924.84"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D37").Value = "3.730.47"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("E39").Value = "  -9.89%  "
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +11.59%  "
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  +2.63%  "

# Numeric-looking text values (keep them stored as text, matching the source)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.23"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.57"
$ws.Range("D6").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.63"
$ws.Range("D14").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.29"
$ws.Range("D18").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.99"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.55"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.704"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.94"
$ws.Range("D24").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.20"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("D28").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("D31").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.10"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.29"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.140"
$ws.Range("D40").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.84"
$ws.Range("D42").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.48"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.59"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "408.05"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.71"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.54"
$ws.Range("D51").ClearFormats()
